# Issue #45 need to debounce button
# - Mark issue #45 (row 45) as Status = DONE (column C)
# - Add two new issues at the bottom of the log:
#     #49  Priority 3, Type UI,   Name "piW fav icon"
#     #50  Priority 3, Type arch, Name "run on port 80"
# - Leave the active selection on C45 (the cell just edited)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Issues")

# Issue #45 is now done
$ws.Range("C45").Value = "DONE"

# New issue #49
$ws.Range("A49").Value = 49
$ws.Range("B49").Value = 3
$ws.Range("D49").Value = "UI"
$ws.Range("E49").Value = "piW fav icon"

# New issue #50
$ws.Range("A50").Value = 50
$ws.Range("B50").Value = 3
$ws.Range("D50").Value = "arch"
$ws.Range("E50").Value = "run on port 80"

# Match the author's final selection/active cell
$ws.Range("C45").Select()
